$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression
$ws.Range("B2").Value = -0.2411497036241854
$ws.Range("C2").Value = -0.2411497036241856
$ws.Range("D2").Value = -0.2411497036241856

# Row 3 - RandomForestRegressor
$ws.Range("B3").Value = 0.9902059634819493
$ws.Range("C3").Value = 0.9895068860712153
$ws.Range("D3").Value = 0.8276613885180325

# Row 4 - rename GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9881013271109148
$ws.Range("C4").Value = 0.9878673286848111
$ws.Range("D4").Value = 0.795051762525095

# Row 5 - rename AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.8739362026414632
$ws.Range("C5").Value = 0.8680699777993957
$ws.Range("D5").Value = 0.5130799457880881
